# Daily update at 8 AM UTC
# Adds the next day's row (row 94) to the "Wins Over Time" tracking sheet,
# and moves the "last row" date-only formatting down to the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet keeps a special date-only number format on the A-cell of the
# final (most recent) row, while every other date in column A uses a
# date+time format. Since row 93 is no longer the last row, restore it to
# the standard date+time format used by the rest of the column.
$ws.Range("A93").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 94.
$ws.Range("A94").Value = 45833
$ws.Range("B94").Value = 400
$ws.Range("C94").Value = 398
$ws.Range("D94").Value = 406

# Row 94 is now the last row, so it gets the date-only format.
$ws.Range("A94").NumberFormat = "YYYY-MM-DD"
